$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: bump the "current refresh time" (column D) for every data row (2-42)
# from 45949.321608796294 to 45950.305879629632.
$newRefresh = 45950.305879629632
for ($r = 2; $r -le 42; $r++) {
    $ws.Cells.Item($r, 4).Value = $newRefresh
}

# --- Step 2: rewrite the rolling "not charged in a while" list, rows 18-42.
# Each row's station (A), terminal (B) and last-charge-end time (C) are replaced
# with the refreshed snapshot (entries that have since charged drop off, newly
# stale ones are appended, the whole list stays sorted by time).

$ws.Cells.Item(18,1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(18,2).Value = "603号直流"
$ws.Cells.Item(18,3).Value = 45948.051631944443
$ws.Cells.Item(19,1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(19,2).Value = "801号直流"
$ws.Cells.Item(19,3).Value = 45948.362951388888
$ws.Cells.Item(20,1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(20,2).Value = "505号直流"
$ws.Cells.Item(20,3).Value = 45948.410497685189
$ws.Cells.Item(21,1).Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Cells.Item(21,2).Value = "111号直流"
$ws.Cells.Item(21,3).Value = 45948.549988425926
$ws.Cells.Item(22,1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(22,2).Value = "902号直流"
$ws.Cells.Item(22,3).Value = 45948.563935185186
$ws.Cells.Item(23,1).Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Cells.Item(23,2).Value = "011A号直流"
$ws.Cells.Item(23,3).Value = 45948.582905092589
$ws.Cells.Item(24,1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(24,2).Value = "705号直流"
$ws.Cells.Item(24,3).Value = 45948.792719907404
$ws.Cells.Item(25,1).Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Cells.Item(25,2).Value = "001A号直流"
$ws.Cells.Item(25,3).Value = 45949.02648148148
$ws.Cells.Item(26,1).Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Cells.Item(26,2).Value = "004A号直流"
$ws.Cells.Item(26,3).Value = 45949.033379629633
$ws.Cells.Item(27,1).Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Cells.Item(27,2).Value = "104号直流"
$ws.Cells.Item(27,3).Value = 45949.101643518516
$ws.Cells.Item(28,1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(28,2).Value = "802号直流"
$ws.Cells.Item(28,3).Value = 45949.187372685185
$ws.Cells.Item(29,1).Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Cells.Item(29,2).Value = "112号直流"
$ws.Cells.Item(29,3).Value = 45949.544594907406
$ws.Cells.Item(30,1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(30,2).Value = "704号直流"
$ws.Cells.Item(30,3).Value = 45949.55395833333
$ws.Cells.Item(31,1).Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Cells.Item(31,2).Value = "208号直流"
$ws.Cells.Item(31,3).Value = 45949.579710648148
$ws.Cells.Item(32,1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(32,2).Value = "803号直流"
$ws.Cells.Item(32,3).Value = 45949.584976851853
$ws.Cells.Item(33,1).Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Cells.Item(33,2).Value = "110号直流"
$ws.Cells.Item(33,3).Value = 45949.610925925925
$ws.Cells.Item(34,1).Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Cells.Item(34,2).Value = "406号直流"
$ws.Cells.Item(34,3).Value = 45949.623159722221
$ws.Cells.Item(35,1).Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Cells.Item(35,2).Value = "205号直流"
$ws.Cells.Item(35,3).Value = 45949.633668981478
$ws.Cells.Item(36,1).Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Cells.Item(36,2).Value = "105号直流"
$ws.Cells.Item(36,3).Value = 45949.639872685184
$ws.Cells.Item(37,1).Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Cells.Item(37,2).Value = "904号直流"
$ws.Cells.Item(37,3).Value = 45949.642638888887
$ws.Cells.Item(38,1).Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Cells.Item(38,2).Value = "210号直流"
$ws.Cells.Item(38,3).Value = 45949.643368055556
$ws.Cells.Item(39,1).Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Cells.Item(39,2).Value = "204号直流"
$ws.Cells.Item(39,3).Value = 45949.646909722222
$ws.Cells.Item(40,1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(40,2).Value = "702号直流"
$ws.Cells.Item(40,3).Value = 45949.727453703701
$ws.Cells.Item(41,1).Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Cells.Item(41,2).Value = "005A号直流"
$ws.Cells.Item(41,3).Value = 45949.743402777778
$ws.Cells.Item(42,1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(42,2).Value = "805号直流"
$ws.Cells.Item(42,3).Value = 45949.766574074078

# --- Step 3: rows 43-52 no longer have stale entries to report -> blank them out
# (keeping their existing cell formatting/styles).
$ws.Range("A43:E52").ClearContents()

# --- Step 4: the author's selection moved from G20 to E12 before saving.
$ws.Range("E12").Select()

